$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3771839289624588
$ws.Range("C2").Value = 3.334164874200428
$ws.Range("B3").Value = 0.3929637689979964
$ws.Range("C3").Value = 4.399788092463178
